# Auto-generated Excel COM-interop script
# Applies numeric cell updates (and a few cell clears/additions) across all 8 sheets
# to match the target OOXML diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1711.6666
$ws.Range("I34").Value = 1711.6666
$ws.Range("K34").Value = 1711.6666
$ws.Range("M34").Value = -1508.6666
$ws.Range("H36").Value = 1711.6666
$ws.Range("I36").Value = 1711.6666
$ws.Range("K36").Value = 1711.6666
$ws.Range("M36").Value = -996.6666
$ws.Range("H43").Value = 3300.6667
$ws.Range("J43").Value = 3701
$ws.Range("L43").Value = 3701
$ws.Range("N43").Value = -3839
$ws.Range("H51").Value = 9990
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H69").Value = 29974.7
$ws.Range("I69").Value = 7112.5
$ws.Range("J69").Value = 45216.168
$ws.Range("K69").Value = 21337.5
$ws.Range("L69").Value = 135648.504
$ws.Range("M69").Value = -20463.5
$ws.Range("N69").Value = -137396.504
$ws.Range("H72").Value = 29974.7
$ws.Range("I72").Value = 7112.5
$ws.Range("J72").Value = 45216.168
$ws.Range("K72").Value = 64012.5
$ws.Range("L72").Value = 406945.512
$ws.Range("M72").Value = -59644.5
$ws.Range("N72").Value = -415681.512
$ws.Range("H98").Value = 1838.25
$ws.Range("I98").Value = 2125
$ws.Range("K98").Value = 2125
$ws.Range("M98").Value = -627
$ws.Range("H113").Value = 8437.25
$ws.Range("I113").Value = 7566.6665
$ws.Range("J113").Value = 8959.6
$ws.Range("K113").Value = 7566.6665
$ws.Range("L113").Value = 8959.6
$ws.Range("M113").Value = -4312.6665
$ws.Range("N113").Value = -15467.6
$ws.Range("H122").Value = 1838.25
$ws.Range("I122").Value = 2125
$ws.Range("K122").Value = 6375
$ws.Range("M122").Value = -3925
$ws.Range("H132").Value = 10396.403
$ws.Range("I132").Value = 1898.25
$ws.Range("K132").Value = 5694.75
$ws.Range("M132").Value = -3164.75
$ws.Range("H137").Value = 3128.025
$ws.Range("I137").Value = 2610.7354
$ws.Range("K137").Value = 7832.206200000001
$ws.Range("M137").Value = -5282.206200000001
$ws.Range("H141").Value = 3592.7407
$ws.Range("I141").Value = 2351.375
$ws.Range("K141").Value = 7054.125
$ws.Range("M141").Value = -1874.125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H5").Value = 150
$ws.Range("I5").Value = 150
$ws.Range("K5").Value = 150
$ws.Range("M5").Value = -38
$ws.Range("H9").Value = 40000
$ws.Range("I9").Value = 75000
$ws.Range("K9").Value = 75000
$ws.Range("M9").Value = -74830
$ws.Range("H20").Value = 40000
$ws.Range("I20").Value = 75000
$ws.Range("K20").Value = 75000
$ws.Range("M20").Value = -74730
$ws.Range("H32").Value = 21741048
$ws.Range("I32").Value = 23257770
$ws.Range("J32").Value = 1379
$ws.Range("K32").Value = 23257770
$ws.Range("L32").Value = 1379
$ws.Range("M32").Value = -23257483
$ws.Range("N32").Value = -1953
$ws.Range("H45").Value = 2260.353
$ws.Range("I45").Value = 1879.1538
$ws.Range("K45").Value = 1879.1538
$ws.Range("M45").Value = -1502.1538
$ws.Range("H63").Value = 18998
$ws.Range("I63").Value = 17996
$ws.Range("J63").Value = 19332
$ws.Range("K63").Value = 17996
$ws.Range("L63").Value = 19332
$ws.Range("M63").Value = -17310
$ws.Range("N63").Value = -20704
$ws.Range("H66").Value = 18998
$ws.Range("I66").Value = 17996
$ws.Range("J66").Value = 19332
$ws.Range("K66").Value = 89980
$ws.Range("L66").Value = 96660
$ws.Range("M66").Value = -86548
$ws.Range("N66").Value = -103524
$ws.Range("H74").Value = 1625.5
$ws.Range("I74").Value = 1690.2759
$ws.Range("J74").Value = 1249.8
$ws.Range("K74").Value = 1690.2759
$ws.Range("L74").Value = 1249.8
$ws.Range("M74").Value = -816.2759000000001
$ws.Range("N74").Value = -2997.8
$ws.Range("H77").Value = 1625.5
$ws.Range("I77").Value = 1690.2759
$ws.Range("J77").Value = 1249.8
$ws.Range("K77").Value = 8451.379500000001
$ws.Range("L77").Value = 6249
$ws.Range("M77").Value = -4083.379500000001
$ws.Range("N77").Value = -14985
$ws.Range("H122").Value = 2264.6
$ws.Range("I122").Value = 1915.75
$ws.Range("K122").Value = 5747.25
$ws.Range("M122").Value = -3297.25
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2629.0588
$ws.Range("I132").Value = 2629.0588
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7887.176399999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5357.176399999999
$ws.Range("N132").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 150
$ws.Range("K4").Value = 150
$ws.Range("M4").Value = -35
$ws.Range("H105").Value = 1623.5555
$ws.Range("I105").Value = 1686.8572
$ws.Range("J105").Value = 1402
$ws.Range("K105").Value = 1686.8572
$ws.Range("L105").Value = 1402
$ws.Range("M105").Value = 60.14280000000008
$ws.Range("N105").Value = -4896
$ws.Range("H134").Value = 1597.0344
$ws.Range("I134").Value = 1404.0714
$ws.Range("K134").Value = 4212.2142
$ws.Range("M134").Value = -1677.2142

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 137576
$ws.Range("I4").Value = 137576
$ws.Range("K4").Value = 137576
$ws.Range("M4").Value = -137464
$ws.Range("H22").Value = 738.125
$ws.Range("I22").Value = 700.7143
$ws.Range("K22").Value = 700.7143
$ws.Range("M22").Value = -350.7143
$ws.Range("H31").Value = 1494.6757
$ws.Range("I31").Value = 1351.5143
$ws.Range("K31").Value = 1351.5143
$ws.Range("M31").Value = -1056.5143
$ws.Range("H34").Value = 1494.6757
$ws.Range("I34").Value = 1351.5143
$ws.Range("K34").Value = 1351.5143
$ws.Range("M34").Value = -1149.5143
$ws.Range("H58").Value = 2011.3478
$ws.Range("I58").Value = 1593.3529
$ws.Range("K58").Value = 1593.3529
$ws.Range("M58").Value = -1390.3529
$ws.Range("H94").Value = 16421.143
$ws.Range("J94").Value = 2587.2
$ws.Range("L94").Value = 2587.2
$ws.Range("N94").Value = -3489.2
$ws.Range("H122").Value = 934861.9399999999
$ws.Range("I122").Value = 3406527
$ws.Range("J122").Value = 7987.5
$ws.Range("K122").Value = 10219581
$ws.Range("L122").Value = 23962.5
$ws.Range("M122").Value = -10217131
$ws.Range("N122").Value = -28862.5
$ws.Range("H134").Value = 4286.0557
$ws.Range("I134").Value = 3596.1765
$ws.Range("K134").Value = 10788.5295
$ws.Range("M134").Value = -8253.529500000001
$ws.Range("H136").Value = 2011.3478
$ws.Range("I136").Value = 1593.3529
$ws.Range("K136").Value = 4780.0587
$ws.Range("M136").Value = -2230.0587

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H87").Value = 33336110
$ws.Range("I87").Value = 33336110
$ws.Range("K87").Value = 100008330
$ws.Range("M87").Value = -100007082
$ws.Range("H90").Value = 33336110
$ws.Range("I90").Value = 33336110
$ws.Range("K90").Value = 300024990
$ws.Range("M90").Value = -300018750
$ws.Range("H113").Value = 751.5
$ws.Range("J113").Value = 727.93335
$ws.Range("L113").Value = 2183.80005
$ws.Range("N113").Value = -6523.80005
$ws.Range("H114").Value = 724.1429000000001
$ws.Range("J114").Value = 960
$ws.Range("L114").Value = 2880
$ws.Range("N114").Value = -9388
$ws.Range("H121").Value = 1302.091
$ws.Range("I121").Value = 756
$ws.Range("K121").Value = 2268
$ws.Range("M121").Value = -958
$ws.Range("H129").Value = 1482.0834
$ws.Range("I129").Value = 777.1429000000001
$ws.Range("J129").Value = 2469
$ws.Range("K129").Value = 2331.4287
$ws.Range("L129").Value = 7407
$ws.Range("M129").Value = 2668.5713
$ws.Range("N129").Value = -17407
$ws.Range("H131").Value = 2166.742
$ws.Range("I131").Value = 1120.619
$ws.Range("K131").Value = 3361.857
$ws.Range("M131").Value = 1678.143

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 11111195
$ws.Range("I2").Value = 36.142857
$ws.Range("K2").Value = 36.142857
$ws.Range("M2").Value = 76.85714300000001
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H94").Value = 30256.385
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 30256.385
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 30256.385
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -31608.385
$ws.Range("H126").Value = 6433.625
$ws.Range("J126").Value = 6645
$ws.Range("L126").Value = 19935
$ws.Range("N126").Value = -24875
$ws.Range("H132").Value = 2009.45
$ws.Range("I132").Value = 1904.6842
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5714.0526
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3184.0526
$ws.Range("N132").Value = -17060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 830.1070999999999
$ws.Range("I22").Value = 749.5
$ws.Range("J22").Value = 836.3077
$ws.Range("K22").Value = 749.5
$ws.Range("L22").Value = 836.3077
$ws.Range("M22").Value = -454.5
$ws.Range("N22").Value = -1426.3077
$ws.Range("H27").Value = 830.1070999999999
$ws.Range("I27").Value = 749.5
$ws.Range("J27").Value = 836.3077
$ws.Range("K27").Value = 749.5
$ws.Range("L27").Value = 836.3077
$ws.Range("M27").Value = -642.5
$ws.Range("N27").Value = -1050.3077
$ws.Range("H55").Value = 1782.6666
$ws.Range("I55").Value = 2177.7144
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 2177.7144
$ws.Range("L55").Value = 400
$ws.Range("M55").Value = -2004.7144
$ws.Range("N55").Value = -746
$ws.Range("H61").Value = 1838.8334
$ws.Range("I61").Value = 1966.7
$ws.Range("J61").Value = 1199.5
$ws.Range("K61").Value = 1966.7
$ws.Range("L61").Value = 1199.5
$ws.Range("M61").Value = -1764.7
$ws.Range("N61").Value = -1603.5
$ws.Range("H113").Value = 1838.8334
$ws.Range("I113").Value = 1966.7
$ws.Range("J113").Value = 1199.5
$ws.Range("K113").Value = 1966.7
$ws.Range("L113").Value = 1199.5
$ws.Range("M113").Value = 203.3
$ws.Range("N113").Value = -5539.5
$ws.Range("H122").Value = 4819.6665
$ws.Range("I122").Value = 3759.5557
$ws.Range("K122").Value = 11278.6671
$ws.Range("M122").Value = -8828.667099999999
$ws.Range("H132").Value = 2542.1316
$ws.Range("I132").Value = 2203.8667
$ws.Range("J132").Value = 3810.625
$ws.Range("K132").Value = 6611.6001
$ws.Range("L132").Value = 11431.875
$ws.Range("M132").Value = -4081.6001
$ws.Range("N132").Value = -16491.875
$ws.Range("H136").Value = 2709.0857
$ws.Range("I136").Value = 2653.1482
$ws.Range("K136").Value = 7959.444600000001
$ws.Range("M136").Value = -5409.444600000001
$ws.Range("H140").Value = 58109.668
$ws.Range("I140").Value = 40000
$ws.Range("J140").Value = 67164.5
$ws.Range("K140").Value = 40000
$ws.Range("L140").Value = 67164.5
$ws.Range("M140").Value = -34820
$ws.Range("N140").Value = -77524.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 2100
$ws.Range("I9").Value = 2100
$ws.Range("K9").Value = 2100
$ws.Range("M9").Value = -1960
$ws.Range("H40").Value = 30247
$ws.Range("J40").Value = 30247
$ws.Range("L40").Value = 30247
$ws.Range("N40").Value = -30545
$ws.Range("H123").Value = 47500
$ws.Range("I123").Value = 30000
$ws.Range("K123").Value = 30000
$ws.Range("M123").Value = -25100
$ws.Range("H135").Value = 68094
$ws.Range("J135").Value = 68094
$ws.Range("L135").Value = 68094
$ws.Range("N135").Value = -78234
